# New weekly Piña price record for Vega Monumental Concepción: insert a row
# above the current row 178, pushing the existing data (rows 178-188) down
# to rows 179-189, then fill the new row 178 with the latest reading.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(178).EntireRow.Insert()

$ws.Range("A178").Value = 11
$ws.Range("B178").Value = "Vega Monumental Concepción"
$ws.Range("C178").Value = "Bíobío"
$ws.Range("D178").Value = 44769
$ws.Range("E178").Value = 8
$ws.Range("F178").Value = "Fruta"
$ws.Range("G178").Value = 100108
$ws.Range("H178").Value = "Tropicales y subtropicales"
$ws.Range("I178").Value = 100108005
$ws.Range("J178").Value = "Piña"
$ws.Range("K178").Value = "Caramelo"
$ws.Range("L178").Value = "Segunda"
$ws.Range("M178").Value = 220
$ws.Range("N178").Value = 18000
$ws.Range("O178").Value = 19000
$ws.Range("P178").Value = 18545
$ws.Range("Q178").Value = "$/caja 14 unidades"
$ws.Range("R178").Value = "Ecuador"
$ws.Range("S178").Value = 1325
$ws.Range("T178").Value = 14
